# -----------------------------------------------------------------------
# Scheduled-runner update: refresh leve-profit calculations (currentAveragePrice
# driven H:N columns) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Generated from the authoritative cell-level diff; each block below sets the
# cells for one worksheet, row by row, left to right (columns H..N).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ========================= Sheet: ALC =========================
$ws = $wb.Worksheets.Item("ALC")

# -- row 6 --
$ws.Range("H6").Value = 1344.4445
$ws.Range("I6").Value = 366.66666
$ws.Range("J6").Value = 1833.3334
$ws.Range("K6").Value = 1099.99998
$ws.Range("L6").Value = 5500.0002
$ws.Range("M6").Value = -987.9999800000001
$ws.Range("N6").Value = -5724.0002
# -- row 125 --
$ws.Range("H125").Value = 1966.4546
$ws.Range("I125").Value = 1872
$ws.Range("J125").Value = 2079.8
$ws.Range("K125").Value = 16848
$ws.Range("L125").Value = 18718.2
$ws.Range("M125").Value = -14388
$ws.Range("N125").Value = -23638.2
# -- row 126 --
$ws.Range("H126").Value = 49900
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 49900
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 49900
$ws.Range("N126").Value = -59780
# -- row 127 --
$ws.Range("H127").Value = 1279.4
$ws.Range("I127").Value = 569.4
$ws.Range("J127").Value = 2699.4
$ws.Range("K127").Value = 1708.2
$ws.Range("L127").Value = 8098.200000000001
$ws.Range("M127").Value = 3251.8
$ws.Range("N127").Value = -18018.2
# -- row 128 --
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
# -- row 129 --
$ws.Range("H129").Value = 5103670.5
$ws.Range("I129").Value = 41667090
$ws.Range("J129").Value = 1798
$ws.Range("K129").Value = 125001270
$ws.Range("L129").Value = 5394
$ws.Range("M129").Value = -124996270
$ws.Range("N129").Value = -15394
# -- row 130 --
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
# -- row 131 --
$ws.Range("H131").Value = 3164.25
$ws.Range("I131").Value = 2627
$ws.Range("J131").Value = 4162
$ws.Range("K131").Value = 7881
$ws.Range("L131").Value = 12486
$ws.Range("M131").Value = -2841
$ws.Range("N131").Value = -22566
# -- row 132 --
$ws.Range("H132").Value = 4083099.2
$ws.Range("I132").Value = 4652113
$ws.Range("J132").Value = 5166.6665
$ws.Range("K132").Value = 13956339
$ws.Range("L132").Value = 15499.9995
$ws.Range("M132").Value = -13953809
$ws.Range("N132").Value = -20559.9995
# -- row 133 --
$ws.Range("H133").Value = 30000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 30000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -40120
# -- row 134 --
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
# -- row 135 --
$ws.Range("H135").Value = 798.1053000000001
$ws.Range("I135").Value = 829.3333
$ws.Range("J135").Value = 236
$ws.Range("K135").Value = 7463.9997
$ws.Range("L135").Value = 2124
$ws.Range("M135").Value = -4928.9997
$ws.Range("N135").Value = -7194
# -- row 136 --
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
# -- row 137 --
$ws.Range("H137").Value = 6182.6665
$ws.Range("I137").Value = 6899.125
$ws.Range("J137").Value = 4749.75
$ws.Range("K137").Value = 20697.375
$ws.Range("L137").Value = 14249.25
$ws.Range("M137").Value = -18147.375
$ws.Range("N137").Value = -19349.25
# -- row 138 --
$ws.Range("H138").Value = 5366.263
$ws.Range("I138").Value = 2274.9524
$ws.Range("J138").Value = 9184.941000000001
$ws.Range("K138").Value = 6824.8572
$ws.Range("L138").Value = 27554.823
$ws.Range("M138").Value = -1684.8572
$ws.Range("N138").Value = -37834.823
# -- row 139 --
$ws.Range("H139").Value = 30000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280
# -- row 140 --
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
# -- row 141 --
$ws.Range("H141").Value = 954837.4
$ws.Range("I141").Value = 2197.1428
$ws.Range("J141").Value = 1907477.6
$ws.Range("K141").Value = 6591.428400000001
$ws.Range("L141").Value = 5722432.800000001
$ws.Range("M141").Value = -1411.428400000001
$ws.Range("N141").Value = -5732792.800000001

# ========================= Sheet: ARM =========================
$ws = $wb.Worksheets.Item("ARM")

# -- row 32 --
$ws.Range("H32").Value = 3889.9028
$ws.Range("I32").Value = 2969.6333
$ws.Range("K32").Value = 2969.6333
$ws.Range("M32").Value = -2682.6333
# -- row 61 --
$ws.Range("H61").Value = 3174.125
$ws.Range("I61").Value = 1596.5
$ws.Range("K61").Value = 1596.5
$ws.Range("M61").Value = -1384.5
# -- row 74 --
$ws.Range("H74").Value = 961.4400000000001
$ws.Range("I74").Value = 850.2857
$ws.Range("J74").Value = 1102.909
$ws.Range("K74").Value = 850.2857
$ws.Range("L74").Value = 1102.909
$ws.Range("M74").Value = 23.71429999999998
$ws.Range("N74").Value = -2850.909
# -- row 77 --
$ws.Range("H77").Value = 961.4400000000001
$ws.Range("I77").Value = 850.2857
$ws.Range("J77").Value = 1102.909
$ws.Range("K77").Value = 4251.4285
$ws.Range("L77").Value = 5514.545
$ws.Range("M77").Value = 116.5715
$ws.Range("N77").Value = -14250.545
# -- row 136 --
$ws.Range("H136").Value = 3174.125
$ws.Range("I136").Value = 1596.5
$ws.Range("K136").Value = 4789.5
$ws.Range("M136").Value = -2239.5

# ========================= Sheet: BSM =========================
$ws = $wb.Worksheets.Item("BSM")

# -- row 20 --
$ws.Range("H20").Value = 4636
$ws.Range("I20").Value = 2440
$ws.Range("J20").Value = 6100
$ws.Range("K20").Value = 2440
$ws.Range("L20").Value = 6100
$ws.Range("M20").Value = -2193
$ws.Range("N20").Value = -6594
# -- row 56 --
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

# ========================= Sheet: CRP =========================
$ws = $wb.Worksheets.Item("CRP")

# -- row 132 --
$ws.Range("H132").Value = 2973.5588
$ws.Range("I132").Value = 2081.9524
$ws.Range("J132").Value = 4413.846
$ws.Range("K132").Value = 6245.8572
$ws.Range("L132").Value = 13241.538
$ws.Range("M132").Value = -3715.8572
$ws.Range("N132").Value = -18301.538

# ========================= Sheet: CUL =========================
$ws = $wb.Worksheets.Item("CUL")

# -- row 131 --
$ws.Range("H131").Value = 1595.9375
$ws.Range("I131").Value = 2043.5625
$ws.Range("J131").Value = 1148.3125
$ws.Range("K131").Value = 6130.6875
$ws.Range("L131").Value = 3444.9375
$ws.Range("M131").Value = -1090.6875
$ws.Range("N131").Value = -13524.9375

# ========================= Sheet: GSM =========================
$ws = $wb.Worksheets.Item("GSM")

# -- row 47 --
$ws.Range("H47").Value = 53773.25
$ws.Range("J47").Value = 53773.25
$ws.Range("L47").Value = 53773.25
$ws.Range("N47").Value = -54909.25
# -- row 124 --
$ws.Range("H124").Value = 19500
$ws.Range("J124").Value = 19500
$ws.Range("L124").Value = 19500
$ws.Range("N124").Value = -29320
# -- row 126 --
$ws.Range("H126").Value = 2918.3076
$ws.Range("I126").Value = 1609.1111
$ws.Range("J126").Value = 3611.4119
$ws.Range("K126").Value = 4827.3333
$ws.Range("L126").Value = 10834.2357
$ws.Range("M126").Value = -2357.3333
$ws.Range("N126").Value = -15774.2357

# ========================= Sheet: LTW =========================
$ws = $wb.Worksheets.Item("LTW")

# -- row 40 --
$ws.Range("H40").Value = 2887.75
$ws.Range("I40").Value = 1899
$ws.Range("J40").Value = 3217.3333
$ws.Range("K40").Value = 1899
$ws.Range("L40").Value = 3217.3333
$ws.Range("M40").Value = -1763
$ws.Range("N40").Value = -3489.3333

# ========================= Sheet: WVR =========================
$ws = $wb.Worksheets.Item("WVR")

# -- row 100 --
$ws.Range("H100").Value = 600.6
$ws.Range("I100").Value = 701.5
$ws.Range("J100").Value = 533.3333
$ws.Range("K100").Value = 1403
$ws.Range("L100").Value = 1066.6666
$ws.Range("M100").Value = -862
$ws.Range("N100").Value = -2148.6666
# -- row 122 --
$ws.Range("H122").Value = 502863.5
$ws.Range("I122").Value = 771343.9399999999
$ws.Range("J122").Value = 4257
$ws.Range("K122").Value = 2314031.82
$ws.Range("L122").Value = 12771
$ws.Range("M122").Value = -2311581.82
$ws.Range("N122").Value = -17671

